$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.78"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.23"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.748"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05802"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.416"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.470"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.318"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1457"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07684"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03248"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03000"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09238"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001679"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.259"
$ws.Range("E16").Value = "15MCDexMCBWorstin24h"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04761"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005994"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006220"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005378"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.690"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1241"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006732"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04289"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007118"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1054"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009724"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005623"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7858"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.09933"
